$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.966.27'
$ws.Range("E2").Value = '  -1.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.907.26'
$ws.Range("E3").Value = '  -3.97%  '

$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.64'
$ws.Range("E5").Value = '  -0.85%  '

$ws.Range("E6").Value = '  +0.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4595'
$ws.Range("E7").Value = '  -1.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3820'
$ws.Range("E8").Value = '  -2.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.54'
$ws.Range("E9").Value = '  -1.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07730'
$ws.Range("E10").Value = '  -2.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9836'
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.10'
$ws.Range("E12").Value = '  -3.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.919.44'
$ws.Range("E13").Value = '  -1.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.993'
$ws.Range("E14").Value = '  -3.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.684'
$ws.Range("E15").Value = '  -3.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07039'
$ws.Range("E16").Value = '  -1.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '84.24'
$ws.Range("E18").Value = '  -4.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009539'
$ws.Range("E19").Value = '  -4.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.72'
$ws.Range("E20").Value = '  -3.50%  '

$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.969.46'
$ws.Range("E22").Value = '  -2.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.330'
$ws.Range("E23").Value = '  -3.47%  '

$ws.Range("E24").Value = '  -2.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.157.79'
$ws.Range("E25").Value = '  -2.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.077'
$ws.Range("E26").Value = '  -0.86%  '

$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.19'
$ws.Range("E28").Value = '  -2.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.592'
$ws.Range("E29").Value = '  -6.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.69'
$ws.Range("E30").Value = '  -2.06%  '

$ws.Range("E31").Value = '  -6.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09261'
$ws.Range("E32").Value = '  -1.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8624'
$ws.Range("E33").Value = '  -4.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.112'
$ws.Range("E34").Value = '  -2.86%  '

$ws.Range("E35").Value = '  -6.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.015'
$ws.Range("E36").Value = '  -5.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05723'
$ws.Range("E37").Value = '  -1.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.147'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.002'
$ws.Range("E39").Value = '  +0.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02044'
$ws.Range("E40").Value = '  -3.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.506'
$ws.Range("E41").Value = '  -4.85%  '

$ws.Range("E42").Value = '  -3.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1756'
$ws.Range("E43").Value = '  -3.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.301'
$ws.Range("E44").Value = '  -5.11%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.745'
$ws.Range("E45").Value = '  +2.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5218'
$ws.Range("E46").Value = '  -2.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.28'
$ws.Range("E47").Value = '  -5.97%  '

$ws.Range("E48").Value = '  -4.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06821'
$ws.Range("E49").Value = '  -1.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.90'
$ws.Range("E50").Value = '  -2.02%  '

$ws.Range("B51").Value = 'PEPE'
$ws.Range("C51").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000002600'
$ws.Range("E51").Value = '  -23.44%  '
